# Edit script: fill in Full Name / Student ID table cells, and
# append the numbered multi-level "Final design" outline list.
$d = $word.ActiveDocument

# --- Table cells: Full Name / Student ID ---
$t = $d.Tables(1)
$t.Cell(1,2).Range.Text = "Hazel Osborne"
$t.Cell(2,2).Range.Text = "1930197"

# --- Build the numbered outline list after the "Final design" title ---
# Locate the "Final design" paragraph (title) to anchor the insertion.
$titlePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "Final design*") {
        $titlePara = $d.Paragraphs($i)
        break
    }
}

$anchorRange = $titlePara.Range
$anchorRange.Collapse(0)
$anchorRange.InsertParagraphAfter()
$newPara = $titlePara.Next()

$first = $true
$prevLevel = 0
$curPara = $newPara

# Item 0: level 0
$curPara.Range.Text = "Output “Welcome to Quest for the Perfect Ice Cream!”"
$curPara.Style = "List Paragraph"
$curPara.Range.ListFormat.ApplyNumberDefault()
$insRange = $curPara.Range
$insRange.Collapse(0)
$insRange.InsertParagraphAfter()
$curPara = $curPara.Next()

# Item 1: level 0
$curPara.Range.Text = "Prompt user to enter their name. “Enter your name”"
$insRange = $curPara.Range
$insRange.Collapse(0)
$insRange.InsertParagraphAfter()
$curPara = $curPara.Next()

# Item 2: level 0
$curPara.Range.Text = "Prompt user to enter their favorite ice cream flavor. “What is your favorite ice cream flavor?”"
$insRange = $curPara.Range
$insRange.Collapse(0)
$insRange.InsertParagraphAfter()
$curPara = $curPara.Next()

# Item 3: level 0
$curPara.Range.Text = "Prompt user to enter a number for variable q1“You’re in the mall and start craving an ice cream. Where do you go? Enter a number below 5 for the hallway, Enter a number above 5 for outside.”"
$insRange = $curPara.Range
$insRange.Collapse(0)
$insRange.InsertParagraphAfter()
$curPara = $curPara.Next()

# Item 4: level 0
$curPara.Range.Text = "Check if q1 < 5 "
$insRange = $curPara.Range
$insRange.Collapse(0)
$insRange.InsertParagraphAfter()
$curPara = $curPara.Next()

# Item 5: level 1
$curPara.Range.Text = "If true, prompt user to enter 1, 2, or 3 for variable q2a. “You see 2 doors. Enter 1 to go into door 1, enter 2 to go to door 2, enter 3 to go back to the main room.” "
$curPara.Range.ListFormat.ListIndent()
$insRange = $curPara.Range
$insRange.Collapse(0)
$insRange.InsertParagraphAfter()
$curPara = $curPara.Next()

# Item 6: level 1
$curPara.Range.Text = "Check if q2a = 1 "
$insRange = $curPara.Range
$insRange.Collapse(0)
$insRange.InsertParagraphAfter()
$curPara = $curPara.Next()

# Item 7: level 2
$curPara.Range.Text = "If true, prompt user to enter polite or rude for q3. “You found an ice cream shop!! Do you want to be rude or polite?”"
$curPara.Range.ListFormat.ListIndent()
$insRange = $curPara.Range
$insRange.Collapse(0)
$insRange.InsertParagraphAfter()
$curPara = $curPara.Next()

# Item 8: level 2
$curPara.Range.Text = "Check if q3 = polite "
$insRange = $curPara.Range
$insRange.Collapse(0)
$insRange.InsertParagraphAfter()
$curPara = $curPara.Next()

# Item 9: level 4
$curPara.Range.Text = "If true, output to user “Congrats!! You got ___ ice cream! Your favorite!”"
$curPara.Range.ListFormat.ListIndent()
$curPara.Range.ListFormat.ListIndent()
$insRange = $curPara.Range
$insRange.Collapse(0)
$insRange.InsertParagraphAfter()
$curPara = $curPara.Next()

# Item 10: level 2
$curPara.Range.Text = "Otherwise, output to user “You got Ice Cream! But because you were rude to the worker, they gave you ____ flavor instead!!"
$curPara.Range.ListFormat.ListOutdent()
$curPara.Range.ListFormat.ListOutdent()
$insRange = $curPara.Range
$insRange.Collapse(0)
$insRange.InsertParagraphAfter()
$curPara = $curPara.Next()

# Item 11: level 1
$curPara.Range.Text = "Otherwise, check if q2a = 2"
$curPara.Range.ListFormat.ListOutdent()
$insRange = $curPara.Range
$insRange.Collapse(0)
$insRange.InsertParagraphAfter()
$curPara = $curPara.Next()

# Item 12: level 2
$curPara.Range.Text = "If true, output to user “You found an ice cream shop!! You ordered your favorite flavor but when you got it, you dropped it and slipped! :( Try again!!” "
$curPara.Range.ListFormat.ListIndent()
$insRange = $curPara.Range
$insRange.Collapse(0)
$insRange.InsertParagraphAfter()
$curPara = $curPara.Next()

# Item 13: level 1
$curPara.Range.Text = "Otherwise, output to user “You leave the hallway and fail to get ice cream! :( Try again!!” "
$curPara.Range.ListFormat.ListOutdent()
$insRange = $curPara.Range
$insRange.Collapse(0)
$insRange.InsertParagraphAfter()
$curPara = $curPara.Next()

# Item 14: level 0
$curPara.Range.Text = "Otherwise, Check if q1>5. "
$curPara.Range.ListFormat.ListOutdent()
$insRange = $curPara.Range
$insRange.Collapse(0)
$insRange.InsertParagraphAfter()
$curPara = $curPara.Next()

# Item 15: level 1
$curPara.Range.Text = "If true, prompt user to enter alley or stall for variable q2b. “You walk outside! You can either go to an unmarked stall, or an alley way. Type ‘alley’ to go to the alleyway or ‘stall’ to go to the stall” "
$curPara.Range.ListFormat.ListIndent()
$insRange = $curPara.Range
$insRange.Collapse(0)
$insRange.InsertParagraphAfter()
$curPara = $curPara.Next()

# Item 16: level 1
$curPara.Range.Text = "Check if q2b = alley "
$insRange = $curPara.Range
$insRange.Collapse(0)
$insRange.InsertParagraphAfter()
$curPara = $curPara.Next()

# Item 17: level 2
$curPara.Range.Text = "If true, output to user “You went into the alley and got scared by a gang of cats! You were so shaken; you went home without ice cream. :( Try again!!”"
$curPara.Range.ListFormat.ListIndent()
$insRange = $curPara.Range
$insRange.Collapse(0)
$insRange.InsertParagraphAfter()
$curPara = $curPara.Next()

# Item 18: level 1
$curPara.Range.Text = "Otherwise, Output to user, “You went to the unmarked stall, and it was an ice cream stall!! You politely ask the worker for your favorite flavor, and they give it to you!! Yay!!”"
$curPara.Range.ListFormat.ListOutdent()
$insRange = $curPara.Range
$insRange.Collapse(0)
$insRange.InsertParagraphAfter()
$curPara = $curPara.Next()

# Item 19: level 0
$curPara.Range.Text = "Otherwise, output to user “Don’t try to cheat the system!! Pick a number below or above five. If not, no ice cream for you!!”"
$curPara.Range.ListFormat.ListOutdent()
$insRange = $curPara.Range
$insRange.Collapse(0)
$insRange.InsertParagraphAfter()
$curPara = $curPara.Next()

# Item 20: level 0
$curPara.Range.Text = "Output to user, “Thanks for playing!”"

# --- Customize the list template levels to match the classic
#     "1. / a. / i." outline numbering pattern (decimal / lowerLetter / lowerRoman) ---
$tmpl = $newPara.Range.ListFormat.ListTemplate
$tmpl.ListLevels(2).NumberStyle = 4  # lowerLetter
$tmpl.ListLevels(3).NumberStyle = 2  # lowerRoman
$tmpl.ListLevels(5).NumberStyle = 4  # lowerLetter
$tmpl.ListLevels(6).NumberStyle = 2  # lowerRoman
$tmpl.ListLevels(8).NumberStyle = 4  # lowerLetter
$tmpl.ListLevels(9).NumberStyle = 2  # lowerRoman

Write-Output "Edit complete."
